$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (header row): shift existing headers one column to the right (B3:J3),
# clear the old A3 value.
$ws.Range("A3").Value = $null
$ws.Range("B3").Value = "split_comp"
$ws.Range("C3").Value = "split"
$ws.Range("D3").Value = "sample"
$ws.Range("E3").Value = "min cluster"
$ws.Range("F3").Value = "n components"
$ws.Range("G3").Value = "training time"
$ws.Range("H3").Value = "nr clusters"
$ws.Range("I3").Value = "topics produced"
$ws.Range("J3").Value = "topic quality (eigene Beurteilung)"

# Row 4: shift one column to the right and set new column A value ("regular")
$ws.Range("A4").Value = "regular"
$ws.Range("B4").Value = "full random"
$ws.Range("C4").Value = "sentences"
$ws.Range("D4").Value = 3000
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = "2000s (ca)"
$ws.Range("H4").Value = "auto"
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = "medium"

# Row 5: shift one column to the right and set new column A value ("regular")
$ws.Range("A5").Value = "regular"
$ws.Range("B5").Value = "full random"
$ws.Range("C5").Value = "sentences"
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = 200
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = "2000s (ca)"
$ws.Range("H5").Value = "auto"
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = "bad"

# Row 6: brand new row for the zeroshot run
$ws.Range("A6").Value = "zeroshot"
$ws.Range("B6").Value = "full random"
$ws.Range("C6").Value = "sentences"
$ws.Range("D6").Value = 5000
$ws.Range("E6").Value = 200
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = "15000s"
$ws.Range("H6").Value = "auto!?"
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = "bad"

# Restore the previously selected cell / cursor position
$ws.Range("E12").Select()
